# Weekly fruit/vegetable price update: rows 2-5 (Caqui, Agrícola del Norte
# S.A. de Arica) are reordered - the data that used to live in row 5 now
# appears in row 2, row 4 moves to row 3, row 3 moves to row 4, and row 2
# moves to row 5 (i.e. rows 2..5 are reversed). Columns that hold identical
# values across all four rows (A,B,C,E,F,G,H,I,J,M,Q,T) show no visible
# change; only D,K,L,N,O,P,R,S differ per diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- old row 5 values
$ws.Range("D2").Value = 44305
$ws.Range("K2").Value = "Mankaki"
$ws.Range("L2").Value = "Segunda"
$ws.Range("N2").Value = 24000
$ws.Range("O2").Value = 25000
$ws.Range("P2").Value = 24500
$ws.Range("R2").Value = "Región de O'Higgins"
$ws.Range("S2").Value = 1361

# Row 3 <- old row 4 values
$ws.Range("D3").Value = 44355
$ws.Range("K3").Value = "Mankaki"
$ws.Range("L3").Value = "Segunda"
$ws.Range("N3").Value = 20000
$ws.Range("O3").Value = 21000
$ws.Range("P3").Value = 20500
$ws.Range("R3").Value = "Región Metropolitana"
$ws.Range("S3").Value = 1139

# Row 4 <- old row 3 values
$ws.Range("D4").Value = 44313
$ws.Range("K4").Value = "Mankaki"
$ws.Range("L4").Value = "Primera"
$ws.Range("N4").Value = 21000
$ws.Range("O4").Value = 22000
$ws.Range("P4").Value = 21500
$ws.Range("R4").Value = "Región de O'Higgins"
$ws.Range("S4").Value = 1194

# Row 5 <- old row 2 values
$ws.Range("D5").Value = 44301
$ws.Range("K5").Value = "Hachiya"
$ws.Range("L5").Value = "Segunda"
$ws.Range("N5").Value = 20000
$ws.Range("O5").Value = 21000
$ws.Range("P5").Value = 20500
$ws.Range("R5").Value = "Región de O'Higgins"
$ws.Range("S5").Value = 1139
